# Daily attendance processing - 2026-01-20 15:16:11
# Normalize the "Recorded By" (column G) values so that the combined
# author list reads "System, dnasr281@gmail.com" instead of
# "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
        $changed++
    }
}

Write-Host "Updated $changed 'Recorded By' cell(s) in column G."
